$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps data between row 2 and row 3 for columns:
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen), S (Precio $/Kg)
$columns = @("D", "M", "N", "O", "P", "R", "S")

foreach ($col in $columns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $val2 = $cell2.Value2
    $val3 = $cell3.Value2
    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
